# Updates cryptos list prices / 1h volume percentages (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 / 17 swapped places (ShibaInu now ranked above WrappedEther) plus refreshed values.
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000136"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.92%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.654.18"
$ws.Range("E17").Value = "  -6.69%  "

# Price (D) / Volume(1h) (E) refresh for every other row.
$ws.Range("D2").Value = "58.849.66"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "2.657.43"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "3.123.31"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "58.827.40"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "339.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.420"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "0.0₃0801"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.911"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.864"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.612"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "275.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0968"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("D47").Value = "2.056.57"
$ws.Range("E47").Value = "  -4.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0533"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0229"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.90%  "
